$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed numeric values (rows 2-6, years 2014/12-2018/12) ---

# Row 2
$ws.Range("D2").Value = 5585
$ws.Range("E2").Value = -37
$ws.Range("F2").Value = -37
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 26
$ws.Range("I2").Value = 26
$ws.Range("K2").Value = 9144
$ws.Range("L2").Value = 2601
$ws.Range("M2").Value = 6543
$ws.Range("N2").Value = 6543
$ws.Range("P2").Value = 837
$ws.Range("Q2").Value = 511
$ws.Range("R2").Value = -323
$ws.Range("S2").Value = -188
$ws.Range("T2").Value = 104
$ws.Range("U2").Value = 406
$ws.Range("V2").Value = 78
$ws.Range("W2").Value = -0.66
$ws.Range("X2").Value = 0.46
$ws.Range("Y2").Value = 0.39
$ws.Range("Z2").Value = 0.28
$ws.Range("AA2").Value = 39.76
$ws.Range("AB2").Value = 715.11
$ws.Range("AC2").Value = 78
$ws.Range("AD2").Value = 106.74
$ws.Range("AE2").Value = 21068
$ws.Range("AF2").Value = 0.4
$ws.Range("AG2").Value = 98
$ws.Range("AH2").Value = 1.17
$ws.Range("AI2").Value = 116.72
$ws.Range("AJ2").Value = 33252697

# Row 3
$ws.Range("D3").Value = 5044
$ws.Range("E3").Value = 385
$ws.Range("F3").Value = 385
$ws.Range("G3").Value = 402
$ws.Range("H3").Value = 326
$ws.Range("I3").Value = 326
$ws.Range("K3").Value = 9217
$ws.Range("L3").Value = 2441
$ws.Range("M3").Value = 6776
$ws.Range("N3").Value = 6776
$ws.Range("P3").Value = 837
$ws.Range("Q3").Value = 490
$ws.Range("R3").Value = -63
$ws.Range("S3").Value = -36
$ws.Range("T3").Value = 125
$ws.Range("U3").Value = 365
$ws.Range("V3").Value = 71
$ws.Range("W3").Value = 7.63
$ws.Range("X3").Value = 6.47
$ws.Range("Y3").Value = 4.9
$ws.Range("Z3").Value = 3.55
$ws.Range("AA3").Value = 36.02
$ws.Range("AB3").Value = 749.73
$ws.Range("AC3").Value = 981
$ws.Range("AD3").Value = 13.65
$ws.Range("AE3").Value = 21818
$ws.Range("AF3").Value = 0.61
$ws.Range("AG3").Value = 195
$ws.Range("AH3").Value = 1.46
$ws.Range("AI3").Value = 18.56
$ws.Range("AJ3").Value = 33252697

# Row 4
$ws.Range("D4").Value = 4680
$ws.Range("E4").Value = -14
$ws.Range("F4").Value = -14
$ws.Range("G4").Value = -232
$ws.Range("H4").Value = -244
$ws.Range("I4").Value = -244
$ws.Range("K4").Value = 8905
$ws.Range("L4").Value = 2370
$ws.Range("M4").Value = 6535
$ws.Range("N4").Value = 6535
$ws.Range("P4").Value = 837
$ws.Range("Q4").Value = 445
$ws.Range("R4").Value = -98
$ws.Range("S4").Value = -75
$ws.Range("T4").Value = 87
$ws.Range("U4").Value = 358
$ws.Range("V4").Value = 68
$ws.Range("W4").Value = -0.31
$ws.Range("X4").Value = -5.22
$ws.Range("Y4").Value = -3.67
$ws.Range("Z4").Value = -2.7
$ws.Range("AA4").Value = 36.26
$ws.Range("AB4").Value = 713.24
$ws.Range("AC4").Value = -735
$ws.Range("AD4").Value = -12.37
$ws.Range("AE4").Value = 21044
$ws.Range("AF4").Value = 0.43
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 33252697

# Row 5
$ws.Range("D5").Value = 4363
$ws.Range("E5").Value = -375
$ws.Range("F5").Value = -375
$ws.Range("G5").Value = -398
$ws.Range("H5").Value = -309
$ws.Range("I5").Value = -309
$ws.Range("K5").Value = 8384
$ws.Range("L5").Value = 2339
$ws.Range("M5").Value = 6045
$ws.Range("N5").Value = 6045
$ws.Range("P5").Value = 837
$ws.Range("Q5").Value = 334
$ws.Range("R5").Value = -479
$ws.Range("S5").Value = -157
$ws.Range("T5").Value = 34
$ws.Range("U5").Value = 300
$ws.Range("V5").Value = 65
$ws.Range("W5").Value = -8.59
$ws.Range("X5").Value = -7.09
$ws.Range("Y5").Value = -4.92
$ws.Range("Z5").Value = -3.58
$ws.Range("AA5").Value = 38.69
$ws.Range("AB5").Value = 674.1900000000001
$ws.Range("AC5").Value = -930
$ws.Range("AD5").Value = -7.86
$ws.Range("AE5").Value = 20722
$ws.Range("AF5").Value = 0.35
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 33252697

# Row 6
$ws.Range("D6").Value = 3963
$ws.Range("E6").Value = 161
$ws.Range("F6").Value = 161
$ws.Range("G6").Value = 260
$ws.Range("H6").Value = 198
$ws.Range("I6").Value = 198
$ws.Range("K6").Value = 8457
$ws.Range("L6").Value = 2382
$ws.Range("M6").Value = 6075
$ws.Range("N6").Value = 6075
$ws.Range("P6").Value = 837
$ws.Range("Q6").Value = 367
$ws.Range("R6").Value = -1285
$ws.Range("S6").Value = -185
$ws.Range("T6").Value = 44
$ws.Range("U6").Value = 323
$ws.Range("V6").Value = 65
$ws.Range("W6").Value = 4.06
$ws.Range("X6").Value = 5
$ws.Range("Y6").Value = 3.27
$ws.Range("Z6").Value = 2.35
$ws.Range("AA6").Value = 39.2
$ws.Range("AB6").Value = 705.27
$ws.Range("AC6").Value = 596
$ws.Range("AD6").Value = 12.08
$ws.Range("AE6").Value = 22707
$ws.Range("AF6").Value = 0.32
$ws.Range("AG6").Value = 195
$ws.Range("AH6").Value = 2.71
$ws.Range("AI6").Value = 26.21
$ws.Range("AJ6").Value = 33252697

# --- Clear cells removed by the edit (J/O columns in rows 2-5) ---
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# --- Clear all forecast data (rows 7-9: 2019(E)-2021(E)), keep only A/B/C ---
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
